# Scheduled market-price refresh for the Tonberry leve-profit tracker.
# For each affected leve row, write the newly observed Universalis average
# prices (currentAveragePrice/NQ/HQ) and the recomputed leve cost/profit
# columns (LevePriceNQ/HQ, LeveProfitNQ/HQ) that derive from them.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart (Roof Tile)
$ws.Range("H19").Value = 1441.1428
$ws.Range("J19").Value = 1964.2222
$ws.Range("L19").Value = 1964.2222
$ws.Range("N19").Value = -2314.2222

# Row 33: Glazed and Confused (Clear Glass Lens)
$ws.Range("H33").Value = 157.36363
$ws.Range("J33").Value = 116.666664
$ws.Range("L33").Value = 116.666664
$ws.Range("N33").Value = -574.666664

# Row 40: Stuck in the Moment (Horn Glue)
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2325

# Row 98: The Dotted Line (Enchanted Durium Ink)
$ws.Range("H98").Value = 985.3333
$ws.Range("I98").Value = 801.28125
$ws.Range("K98").Value = 801.28125
$ws.Range("M98").Value = 696.71875

# Row 106: Making Your Mark (Enchanted Palladium Ink)
$ws.Range("H106").Value = 4186.6
$ws.Range("I106").Value = 4186.6
$ws.Range("K106").Value = 4186.6
$ws.Range("M106").Value = -3555.6

# Row 122: Wishful Inking (Enchanted High Durium Ink)
$ws.Range("H122").Value = 985.3333
$ws.Range("I122").Value = 801.28125
$ws.Range("K122").Value = 2403.84375
$ws.Range("M122").Value = 46.15625

# Row 135: For Tired Minds (Grade 1 Gemsap of Intelligence)
$ws.Range("H135").Value = 52632484
$ws.Range("I135").Value = 1034.3846
$ws.Range("K135").Value = 9309.4614
$ws.Range("M135").Value = -6774.4614

# Row 137: Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws.Range("H137").Value = 2188.4546
$ws.Range("I137").Value = 1737.0769
$ws.Range("K137").Value = 5211.2307
$ws.Range("M137").Value = -2661.2307

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust (Steel Ingot)
$ws.Range("H32").Value = 4987.314
$ws.Range("I32").Value = 3830.0476
$ws.Range("K32").Value = 3830.0476
$ws.Range("M32").Value = -3543.0476

# Row 74: As the Bolt Flies (Titanium Nugget)
$ws.Range("H74").Value = 1123.0834
$ws.Range("I74").Value = 573.8276
$ws.Range("K74").Value = 573.8276
$ws.Range("M74").Value = 300.1724

# Row 77: Heavy Metal Banned (L) (Titanium Nugget)
$ws.Range("H77").Value = 1123.0834
$ws.Range("I77").Value = 573.8276
$ws.Range("K77").Value = 2869.138
$ws.Range("M77").Value = 1498.862

# Row 109: A Head of Demand (Deepgold Helm of Fending)
$ws.Range("H109").Value = 67888
$ws.Range("J109").Value = 67888
$ws.Range("L109").Value = 67888
$ws.Range("N109").Value = -70662

# Row 133: Shielding My Students (Mountain Chromite Tower Shield)
$ws.Range("H133").Value = 6000
$ws.Range("J133").Value = 6000
$ws.Range("L133").Value = 6000
$ws.Range("N133").Value = -11060

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run (Iron Rivets)
$ws.Range("H22").Value = 600.3333
$ws.Range("I22").Value = 600.3333
$ws.Range("K22").Value = 600.3333
$ws.Range("M22").Value = -427.3333

# Row 99: Meddle in Metal (Oroshigane Ingot)
$ws.Range("H99").Value = 1200.25
$ws.Range("I99").Value = 899.2
$ws.Range("K99").Value = 899.2
$ws.Range("M99").Value = 598.8

# Row 107: The Gold Experience (Deepgold Nugget)
$ws.Range("H107").Value = 622.2143
$ws.Range("I107").Value = 501
$ws.Range("K107").Value = 501
$ws.Range("M107").Value = 1419

$ws = $wb.Worksheets.Item("CRP")
# Row 62: Splinter in the Sewers (Cedar Lumber)
$ws.Range("H62").Value = 2314.2856
$ws.Range("I62").Value = 2440
$ws.Range("K62").Value = 2440
$ws.Range("M62").Value = -1816

# Row 65: The Lumber of Their Discontent (L) (Cedar Lumber)
$ws.Range("H65").Value = 2314.2856
$ws.Range("I65").Value = 2440
$ws.Range("K65").Value = 12200
$ws.Range("M65").Value = -9080

# Row 99: O Pine (Pine Lumber)
$ws.Range("H99").Value = 3064.375
$ws.Range("I99").Value = 2073.5715
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 2073.5715
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -575.5715
$ws.Range("N99").Value = -12996

# Row 105: Zelkova, My Love (Zelkova Lumber)
$ws.Range("H105").Value = 1053.1
$ws.Range("I105").Value = 1059
$ws.Range("K105").Value = 1059
$ws.Range("M105").Value = 688

# Row 126: A Better Conductor (Red Pine Lumber)
$ws.Range("H126").Value = 3064.375
$ws.Range("I126").Value = 2073.5715
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 6220.7145
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -3750.7145
$ws.Range("N126").Value = -34940

# Row 132: Hull Lotta Damage (Ginseng Lumber)
$ws.Range("H132").Value = 3289.05
$ws.Range("I132").Value = 2528.9167
$ws.Range("K132").Value = 7586.750100000001
$ws.Range("M132").Value = -5056.750100000001

# Row 134: Wood You Be Quiet (Ceiba Lumber)
$ws.Range("H134").Value = 2714.12
$ws.Range("I134").Value = 2292.85
$ws.Range("K134").Value = 6878.549999999999
$ws.Range("M134").Value = -4343.549999999999

# Row 141: No Greater Treasure (Claro Walnut Necklace of Gathering)
$ws.Range("H141").Value = 72000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food (Table Salt)
$ws.Range("H2").Value = 192.91667
$ws.Range("I2").Value = 155.875
$ws.Range("K2").Value = 935.25
$ws.Range("M2").Value = -822.25

# Row 38: Pretty as a Picture (Dark Vinegar)
$ws.Range("H38").Value = 540.4286
$ws.Range("J38").Value = 1000.6667
$ws.Range("L38").Value = 3002.0001
$ws.Range("N38").Value = -3696.0001

# Row 59: Comfort Me with Mushrooms (Buttons in a Blanket)
$ws.Range("H59").Value = 650
$ws.Range("I59").Value = 475
$ws.Range("J59").Value = 1000
$ws.Range("K59").Value = 1425
$ws.Range("L59").Value = 3000
$ws.Range("M59").Value = -885
$ws.Range("N59").Value = -4080

# Row 122: Salt of the North (Northern Sea Salt)
$ws.Range("H122").Value = 2591.6
$ws.Range("J122").Value = 2591.6
$ws.Range("L122").Value = 23324.4
$ws.Range("N122").Value = -28224.4

# Row 137: Creative Chocolate (Gateau au Chocolat)
$ws.Range("H137").Value = 4203.16
$ws.Range("I137").Value = 1645.3636
$ws.Range("J137").Value = 6212.857
$ws.Range("K137").Value = 4936.0908
$ws.Range("L137").Value = 18638.571
$ws.Range("M137").Value = 163.9092000000001
$ws.Range("N137").Value = -28838.571

$ws = $wb.Worksheets.Item("GSM")
# Row 92: Play It by Ear (Triphane Earrings of Healing)
$ws.Range("H92").Value = 23919.4
$ws.Range("J92").Value = 23919.4
$ws.Range("L92").Value = 23919.4
$ws.Range("N92").Value = -27663.4

# Row 102: Put the Metal to the Peddle (Durium Ingot)
$ws.Range("H102").Value = 3642.724
$ws.Range("I102").Value = 4171.421
$ws.Range("J102").Value = 2638.2
$ws.Range("K102").Value = 4171.421
$ws.Range("L102").Value = 2638.2
$ws.Range("M102").Value = -2549.421
$ws.Range("N102").Value = -5882.2

# Row 113: Copious Crystal Cannons (Manasilver Nugget)
$ws.Range("H113").Value = 1044.9166
$ws.Range("I113").Value = 645.2857
$ws.Range("J113").Value = 1604.4
$ws.Range("K113").Value = 645.2857
$ws.Range("L113").Value = 1604.4
$ws.Range("M113").Value = 1524.7143
$ws.Range("N113").Value = -5944.4

# Row 122: Awarding Academic Excellence (Ametrine)
$ws.Range("H122").Value = 1095.7407
$ws.Range("I122").Value = 1111.88
$ws.Range("K122").Value = 3335.64
$ws.Range("M122").Value = -885.6400000000003

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad (Toad Leather)
$ws.Range("H40").Value = 13858.083
$ws.Range("I40").Value = 15112.875
$ws.Range("J40").Value = 11348.5
$ws.Range("K40").Value = 15112.875
$ws.Range("L40").Value = 11348.5
$ws.Range("M40").Value = -14976.875
$ws.Range("N40").Value = -11620.5

# Row 46: Supply Side Logic (Boar Leather)
$ws.Range("H46").Value = 1414
$ws.Range("I46").Value = 738.7
$ws.Range("J46").Value = 2378.7144
$ws.Range("K46").Value = 738.7
$ws.Range("L46").Value = 2378.7144
$ws.Range("M46").Value = -550.7
$ws.Range("N46").Value = -2754.7144

# Row 55: It's Not a Job, It's a Calling (Peiste Leather)
$ws.Range("H55").Value = 431.8889
$ws.Range("I55").Value = 423.86667
$ws.Range("J55").Value = 441.91666
$ws.Range("K55").Value = 423.86667
$ws.Range("L55").Value = 441.91666
$ws.Range("M55").Value = -250.86667
$ws.Range("N55").Value = -787.91666

$ws = $wb.Worksheets.Item("WVR")
# Row 61: Bundle Up, It's Odd out There (Woolen Deerstalker)
$ws.Range("H61").Value = 9000
$ws.Range("J61").Value = 9000
$ws.Range("L61").Value = 9000
$ws.Range("N61").Value = -9584

# Row 100: Of Great Import (Kudzu Thread)
$ws.Range("H100").Value = 716
$ws.Range("I100").Value = 492
$ws.Range("K100").Value = 984
$ws.Range("M100").Value = -443

# Row 113: A Tender Table (Pixie Floss)
$ws.Range("H113").Value = 445.4091
$ws.Range("I113").Value = 326.41666
$ws.Range("K113").Value = 979.2499799999999
$ws.Range("M113").Value = 1190.75002

# Row 123: Helping Handwear (Fingerless Darkhempen Gloves of Healing)
$ws.Range("H123").Value = 47630.555
$ws.Range("J123").Value = 47630.555
$ws.Range("L123").Value = 47630.555
$ws.Range("N123").Value = -57430.555
